# Update "想去人数" (column F) values across all four sheets to reflect the
# newly scraped counts published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 22
$ws.Range("F3").Value = 355
$ws.Range("F4").Value = 1351
$ws.Range("F7").Value = 3955
$ws.Range("F8").Value = 246
$ws.Range("F9").Value = 804
$ws.Range("F10").Value = 2414
$ws.Range("F11").Value = 377
$ws.Range("F13").Value = 245
$ws.Range("F14").Value = 766
$ws.Range("F15").Value = 222
$ws.Range("F16").Value = 210
$ws.Range("F17").Value = 4117
$ws.Range("F21").Value = 361
$ws.Range("F23").Value = 61

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 60
$ws.Range("F4").Value = 35
$ws.Range("F8").Value = 108
$ws.Range("F12").Value = 10
$ws.Range("F23").Value = 91

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 836
$ws.Range("F4").Value = 2146
$ws.Range("F6").Value = 32
$ws.Range("F7").Value = 2

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 836
$ws.Range("F3").Value = 2146
$ws.Range("F5").Value = 60
$ws.Range("F6").Value = 22
$ws.Range("F8").Value = 355
$ws.Range("F9").Value = 1351
$ws.Range("F11").Value = 35
$ws.Range("F14").Value = 32
$ws.Range("F16").Value = 3955
$ws.Range("F18").Value = 246
$ws.Range("F19").Value = 108
$ws.Range("F21").Value = 804
$ws.Range("F22").Value = 2414
$ws.Range("F23").Value = 377
$ws.Range("F26").Value = 245
$ws.Range("F27").Value = 766
$ws.Range("F28").Value = 222
$ws.Range("F29").Value = 210
$ws.Range("F30").Value = 10
$ws.Range("F40").Value = 61
$ws.Range("F47").Value = 2
$ws.Range("F48").Value = 91
